$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.525.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.836.94"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -2.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4303"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3710"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8689"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.21"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.70"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.707"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.379"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07093"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008943"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.534.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.178"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.068.59"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.005"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.50"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.156"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.308"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08882"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7688"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.499"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.906"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.007"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.126"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01963"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05292"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.167"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.883"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5094"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1677"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.713"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.54"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4728"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("E48").Value = "  -1.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.674"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.834"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.76%  "
